$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the currency name and rate in row 2
$ws.Range("A2").Value = "BGN"
$ws.Range("B2").Value = 44.204898834228516

# Remove row 3 entirely (shifts nothing below it, just clears/removes the row)
$ws.Range("A3:B3").Delete()
